$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.966.28"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").Value = "1.874.67"
$ws.Range("E3").Value = "  +0.50%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'305.46"
$ws.Range("E5").Value = "  -0.27%  "
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").Value = "'0.5087"
$ws.Range("E7").Value = "  -1.09%  "
$ws.Range("D8").Value = "'0.3671"
$ws.Range("E8").Value = "  -2.21%  "
$ws.Range("D9").Value = "'0.07210"
$ws.Range("E9").Value = "  +0.66%  "
$ws.Range("D10").Value = "'0.8956"
$ws.Range("E10").Value = "  +0.77%  "
$ws.Range("D11").Value = "'20.78"
$ws.Range("E11").Value = "  +0.50%  "
$ws.Range("D12").Value = "1.897.66"
$ws.Range("E12").Value = "  +2.07%  "
$ws.Range("D13").Value = "'0.07523"
$ws.Range("E13").Value = "  -0.98%  "
$ws.Range("D14").Value = "'95.16"
$ws.Range("E14").Value = "  +6.37%  "
$ws.Range("D15").Value = "'5.248"
$ws.Range("E15").Value = "  -1.11%  "
$ws.Range("E16").Value = "  +0.00%  "
$ws.Range("D17").Value = "'0.000008537"
$ws.Range("E18").Value = "  +1.36%  "
$ws.Range("D20").Value = "27.005.26"
$ws.Range("E20").Value = "  -0.30%  "
$ws.Range("D21").Value = "'5.026"
$ws.Range("D22").Value = "2.125.82"
$ws.Range("E22").Value = "  +2.45%  "
$ws.Range("D23").Value = "'10.39"
$ws.Range("E23").Value = "  -1.18%  "
$ws.Range("D24").Value = "'6.407"
$ws.Range("E24").Value = "  -0.72%  "
$ws.Range("D25").Value = "'148.58"
$ws.Range("E25").Value = "  +0.82%  "
$ws.Range("E26").Value = "  -2.78%  "
$ws.Range("D27").Value = "'17.91"
$ws.Range("E27").Value = "  -0.42%  "
$ws.Range("D28").Value = "'2.094"
$ws.Range("E28").Value = "  -1.04%  "
$ws.Range("D29").Value = "'113.52"
$ws.Range("E29").Value = "  +0.72%  "
$ws.Range("D30").Value = "'4.729"
$ws.Range("E30").Value = "  +1.44%  "
$ws.Range("D31").Value = "'4.738"
$ws.Range("E31").Value = "  +0.78%  "
$ws.Range("D32").Value = "'0.09159"
$ws.Range("E32").Value = "  +0.63%  "
$ws.Range("D33").Value = "'0.05105"
$ws.Range("E33").Value = "  -0.46%  "
$ws.Range("D34").Value = "'0.7508"
$ws.Range("D35").Value = "'2.968"
$ws.Range("E35").Value = "  -3.27%  "
$ws.Range("D36").Value = "'1.160"
$ws.Range("E36").Value = "  +0.27%  "
$ws.Range("D37").Value = "'3.237"
$ws.Range("E37").Value = "  +6.41%  "
$ws.Range("D38").Value = "'2.536"
$ws.Range("E38").Value = "  +2.06%  "
$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").Value = "'0.5643"
$ws.Range("E39").Value = "  +5.71%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.02004"
$ws.Range("E40").Value = "  -1.89%  "
$ws.Range("D41").Value = "'1.078"
$ws.Range("E41").Value = "  +0.34%  "
$ws.Range("D42").Value = "'6.646"
$ws.Range("E42").Value = "  +1.39%  "
$ws.Range("D43").Value = "'115.38"
$ws.Range("E43").Value = "  -1.19%  "
$ws.Range("D44").Value = "'8.587"
$ws.Range("E44").Value = "  +3.73%  "
$ws.Range("D45").Value = "'0.1480"
$ws.Range("E45").Value = "  +0.67%  "
$ws.Range("D46").Value = "'0.4754"
$ws.Range("E46").Value = "  +2.46%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'10.17"
$ws.Range("E47").Value = "  +1.81%  "
$ws.Range("B48").Value = "PaxDollar"
$ws.Range("C48").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D48").Value = "'0.9998"
$ws.Range("E48").Value = "  -0.07%  "
$ws.Range("D49").Value = "'1.574"
$ws.Range("E49").Value = "  +0.15%  "
$ws.Range("E50").Value = "  +1.17%  "
$ws.Range("D51").Value = "'63.26"
$ws.Range("E51").Value = "  -0.93%  "
